# Auto-generated Excel COM-interop script to update the cryptocurrency price table
# on Sheet1 (columns B=Coin, C=Link, D=Price, E=Volume(1h)).
#
# All D/E values in this sheet are plain text (not real numbers/percentages),
# e.g. '574.55', '3.180.37', '  +2.53%  '. Writing such strings via .Value would
# normally get auto-coerced by Excel into floating point numbers, which would
# corrupt values like '65.129.04' or lose exact text like '1.00' -> 1. To avoid
# that, we snapshot the existing cell style, temporarily force Text number format,
# assign the literal string, then restore the original style so no visible
# formatting/styling changes are introduced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$origStyle = $ws.Range("D2").Style

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = $origStyle
}

# Row 2
Set-TextValue "D2" '65.129.04'
Set-TextValue "E2" '  +1.68%  '

# Row 3
Set-TextValue "D3" '3.180.37'

# Row 4
Set-TextValue "E4" '  -0.03%  '

# Row 5
Set-TextValue "D5" '574.55'
Set-TextValue "E5" '  +2.53%  '

# Row 6
Set-TextValue "D6" '151.17'
Set-TextValue "E6" '  +5.04%  '

# Row 7
Set-TextValue "E7" '  -0.11%  '

# Row 8
Set-TextValue "D8" '3.179.09'

# Row 9
Set-TextValue "E9" '  +3.33%  '

# Row 10
Set-TextValue "E10" '  +4.67%  '

# Row 11
Set-TextValue "D11" '6.23'
Set-TextValue "E11" '  +2.39%  '

# Row 12
Set-TextValue "D12" '0.507'
Set-TextValue "E12" '  +4.86%  '

# Row 13
Set-TextValue "D13" '0.0000275'
Set-TextValue "E13" '  +18.73%  '

# Row 14
Set-TextValue "D14" '38.23'
Set-TextValue "E14" '  +7.68%  '

# Row 15
Set-TextValue "D15" '3.697.53'
Set-TextValue "E15" '  +3.77%  '

# Row 16
Set-TextValue "D16" '65.218.39'
Set-TextValue "E16" '  +1.71%  '

# Row 17
Set-TextValue "B17" 'WrappedEther'
Set-TextValue "C17" 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue "D17" '3.186.15'
Set-TextValue "E17" '  +3.82%  '

# Row 18
Set-TextValue "B18" 'Polkadot'
Set-TextValue "C18" 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue "D18" '7.22'
Set-TextValue "E18" '  +6.70%  '

# Row 19
Set-TextValue "E19" '  +1.20%  '

# Row 20
Set-TextValue "D20" '513.13'
Set-TextValue "E20" '  +7.19%  '

# Row 21
Set-TextValue "D21" '14.97'
Set-TextValue "E21" '  +7.12%  '

# Row 22
Set-TextValue "E22" '  +8.41%  '

# Row 23
Set-TextValue "D23" '15.63'
Set-TextValue "E23" '  +9.00%  '

# Row 24
Set-TextValue "E24" '  +3.71%  '

# Row 25
Set-TextValue "D25" '85.09'
Set-TextValue "E25" '  +3.70%  '

# Row 26
Set-TextValue "E26" '  +0.11%  '

# Row 27
Set-TextValue "D27" '9.18'
Set-TextValue "E27" '  +14.54%  '

# Row 28
Set-TextValue "D28" '2.92'
Set-TextValue "E28" '  +4.12%  '

# Row 29
Set-TextValue "E29" '  +8.28%  '

# Row 30
Set-TextValue "D30" '28.18'
Set-TextValue "E30" '  +6.90%  '

# Row 31
Set-TextValue "D31" '2.80'
Set-TextValue "E31" '  +14.50%  '

# Row 32
Set-TextValue "E32" '  +7.46%  '

# Row 33
Set-TextValue "D33" '1.00'
Set-TextValue "E33" '  -0.02%  '

# Row 34
Set-TextValue "E34" '  +11.15%  '

# Row 35
Set-TextValue "D35" '6.72'
Set-TextValue "E35" '  +7.53%  '

# Row 36
Set-TextValue "D36" '55.77'
Set-TextValue "E36" '  +1.51%  '

# Row 37
Set-TextValue "B37" 'Hedera'
Set-TextValue "C37" 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue "D37" '0.0894'
Set-TextValue "E37" '  +10.10%  '

# Row 38
Set-TextValue "B38" 'Bittensor'
Set-TextValue "C38" 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue "D38" '481.29'
Set-TextValue "E38" '  +7.99%  '

# Row 39
Set-TextValue "D39" '3.10'
Set-TextValue "E39" '  +8.79%  '

# Row 40
Set-TextValue "D40" '0.0422'
Set-TextValue "E40" '  +3.31%  '

# Row 41
Set-TextValue "D41" '3.144.26'
Set-TextValue "E41" '  +4.78%  '

# Row 42
Set-TextValue "D42" '8.66'
Set-TextValue "E42" '  +4.92%  '

# Row 43
Set-TextValue "E43" '  +4.51%  '

# Row 44
Set-TextValue "D44" '2.51'
Set-TextValue "E44" '  +15.93%  '

# Row 45
Set-TextValue "D45" '0.289'
Set-TextValue "E45" '  +10.59%  '

# Row 46
Set-TextValue "D46" '29.37'
Set-TextValue "E46" '  +5.07%  '

# Row 47
Set-TextValue "D47" '0.0₃0601'
Set-TextValue "E47" '  +15.56%  '

# Row 48
Set-TextValue "E48" '  -0.05%  '

# Row 49
Set-TextValue "E49" '  +2.27%  '

# Row 50
Set-TextValue "D50" '2.30'
Set-TextValue "E50" '  +11.16%  '

# Row 51
Set-TextValue "D51" '122.60'
Set-TextValue "E51" '  +3.31%  '
